$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "22.372.31"
$ws.Cells.Item(2, 5).Value = "  -4.42%  "

$ws.Cells.Item(3, 4).Value = "1.567.73"
$ws.Cells.Item(3, 5).Value = "  -4.65%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.002"
$ws.Cells.Item(4, 5).Value = "  -0.01%  "

$ws.Cells.Item(5, 5).Value = "  +0.01%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "290.69"
$ws.Cells.Item(6, 5).Value = "  -2.85%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.3681"
$ws.Cells.Item(7, 5).Value = "  -2.81%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "49.43"

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.3397"
$ws.Cells.Item(9, 5).Value = "  -3.55%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "1.172"
$ws.Cells.Item(10, 5).Value = "  -3.32%  "

$ws.Cells.Item(11, 5).Value = "  -5.90%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.002"
$ws.Cells.Item(12, 5).Value = "  -0.05%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "21.21"
$ws.Cells.Item(13, 5).Value = "  -3.79%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.050"
$ws.Cells.Item(14, 5).Value = "  -5.00%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "6.891"
$ws.Cells.Item(15, 5).Value = "  -5.72%  "

$ws.Cells.Item(16, 2).Value = "WrappedEther"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(16, 4).Value = "1.578.17"
$ws.Cells.Item(16, 5).Value = "  -3.74%  "

$ws.Cells.Item(17, 2).Value = "ShibaInu"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.00001138"
$ws.Cells.Item(17, 5).Value = "  -5.26%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "89.23"
$ws.Cells.Item(18, 5).Value = "  -7.75%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.06776"
$ws.Cells.Item(19, 5).Value = "  -3.01%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "1.002"
$ws.Cells.Item(20, 5).Value = "  +0.09%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.234"
$ws.Cells.Item(21, 5).Value = "  -7.21%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.5357"
$ws.Cells.Item(22, 5).Value = "  -6.10%  "

$ws.Cells.Item(23, 5).Value = "  -4.93%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "12.04"
$ws.Cells.Item(24, 5).Value = "  -2.58%  "

$ws.Cells.Item(25, 4).Value = "22.379.81"
$ws.Cells.Item(25, 5).Value = "  -4.50%  "

$ws.Cells.Item(26, 5).Value = "  -4.26%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "2.976"
$ws.Cells.Item(27, 5).Value = "  +2.61%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "19.89"
$ws.Cells.Item(28, 5).Value = "  -4.52%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "145.70"
$ws.Cells.Item(29, 5).Value = "  -4.42%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "4.954"
$ws.Cells.Item(30, 5).Value = "  -4.86%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "125.39"
$ws.Cells.Item(31, 5).Value = "  -5.41%  "

$ws.Cells.Item(32, 4).Value = "1.758.51"
$ws.Cells.Item(32, 5).Value = "  -3.51%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.038"
$ws.Cells.Item(33, 5).Value = "  +5.78%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "6.248"
$ws.Cells.Item(34, 5).Value = "  -9.05%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.992"
$ws.Cells.Item(35, 5).Value = "  -6.05%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "10.32"
$ws.Cells.Item(36, 5).Value = "  -9.46%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.08461"
$ws.Cells.Item(37, 5).Value = "  -3.24%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.02543"
$ws.Cells.Item(38, 5).Value = "  -5.77%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.2331"
$ws.Cells.Item(39, 5).Value = "  -4.01%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.06552"
$ws.Cells.Item(40, 5).Value = "  -3.59%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "5.535"
$ws.Cells.Item(41, 5).Value = "  -6.25%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "11.83"
$ws.Cells.Item(42, 5).Value = "  -7.70%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "1.246"
$ws.Cells.Item(43, 5).Value = "  -3.45%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.6372"
$ws.Cells.Item(44, 5).Value = "  -7.03%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "14.40"
$ws.Cells.Item(45, 5).Value = "  -7.62%  "

$ws.Cells.Item(46, 5).Value = "  -0.07%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.6008"
$ws.Cells.Item(47, 5).Value = "  -5.03%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "3.778"
$ws.Cells.Item(48, 5).Value = "  -3.28%  "

$ws.Cells.Item(49, 5).Value = "  -5.16%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.267"
$ws.Cells.Item(50, 5).Value = "  +7.76%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "123.95"
$ws.Cells.Item(51, 5).Value = "  -2.56%  "
